$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-27 08:18:25"
$wsZhCn.Range("G3").Value = "2016-01-27 08:19:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-27 08:18:37"
$wsDeDe.Range("G3").Value = "2016-01-27 08:19:31"
